$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = Get-Date -Year 2023 -Month 9 -Day 13 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 289; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
